$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the date value in A2 (45743 -> 45746, i.e. 2025-03-27 -> 2025-03-30)
$ws.Range("A2").Value = 45746

# Update the active selection on the sheet from I10 to D7
$ws.Range("D7").Select()
